# Weekly update: insert a new week's worth of data (2 rows) for
# Femacal de La Calera - Limon, at the top of the price log, pushing
# the previously-most-recent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 527 and 528 (this shifts the existing
# rows 527..643 down to 529..645 and grows the sheet dimension to
# A1:T645). Row formatting (e.g. the date style on column D) carries
# over automatically from the row above, just like a manual Excel
# "Insert Rows" would do.
$ws.Range("A527:A528").EntireRow.Insert()

# The two freshly-inserted rows are completely blank, so every column
# needs to be (re)written, not just the ones that differ from the row
# they displaced. Columns A, B, C, E-K, Q, R and T are identical on
# every single record in this sheet (same market / product / unit),
# so reuse those constants for the two new rows.
$marketId = 3
$market = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad = "Sin especificar"
$unidad = "$/malla 16 kilos"
$origen = "Provincia de Quillota"
$kgUnidad = 16

# Row 527 keeps its quality grade ("1a amarillo") but gets the new
# week's figures.
$ws.Range("A527").Value = $marketId
$ws.Range("B527").Value = $market
$ws.Range("C527").Value = $region
$ws.Range("D527").Value = 44511
$ws.Range("E527").Value = $codreg
$ws.Range("F527").Value = $tipo
$ws.Range("G527").Value = $productoId
$ws.Range("H527").Value = $producto
$ws.Range("I527").Value = $categoriaId
$ws.Range("J527").Value = $categoria
$ws.Range("K527").Value = $variedad
$ws.Range("L527").Value = "1a amarillo"
$ws.Range("M527").Value = 230
$ws.Range("N527").Value = 4500
$ws.Range("O527").Value = 5000
$ws.Range("P527").Value = 4783
$ws.Range("Q527").Value = $unidad
$ws.Range("R527").Value = $origen
$ws.Range("S527").Value = 299
$ws.Range("T527").Value = $kgUnidad

# Row 528 becomes the new week's "2a amarillo" entry (it used to be
# "1a plateado" before the insert).
$ws.Range("A528").Value = $marketId
$ws.Range("B528").Value = $market
$ws.Range("C528").Value = $region
$ws.Range("D528").Value = 44511
$ws.Range("E528").Value = $codreg
$ws.Range("F528").Value = $tipo
$ws.Range("G528").Value = $productoId
$ws.Range("H528").Value = $producto
$ws.Range("I528").Value = $categoriaId
$ws.Range("J528").Value = $categoria
$ws.Range("K528").Value = $variedad
$ws.Range("L528").Value = "2a amarillo"
$ws.Range("M528").Value = 308
$ws.Range("N528").Value = 3000
$ws.Range("O528").Value = 4000
$ws.Range("P528").Value = 3565
$ws.Range("Q528").Value = $unidad
$ws.Range("R528").Value = $origen
$ws.Range("S528").Value = 223
$ws.Range("T528").Value = $kgUnidad
